$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D to remain text so numeric-looking strings (e.g. "1.00", "3.00")
# are not coerced into numbers and lose their formatting / trailing zeros.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "70.926.53"
$ws.Range("E2").Value = "  +1.79%  "
$ws.Range("D3").Value = "3.637.87"
$ws.Range("E3").Value = "  +3.76%  "
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  +0.10%  "
$ws.Range("D5").Value = "606.09"
$ws.Range("E5").Value = "  +0.06%  "
$ws.Range("D6").Value = "200.14"
$ws.Range("E6").Value = "  +2.45%  "
$ws.Range("D7").Value = "0.629"
$ws.Range("E7").Value = "  +0.30%  "
$ws.Range("D8").Value = "1.00"
$ws.Range("E8").Value = "  +0.09%  "
$ws.Range("D9").Value = "0.220"
$ws.Range("E9").Value = "  +8.90%  "
$ws.Range("D10").Value = "0.649"
$ws.Range("E10").Value = "  -0.15%  "
$ws.Range("D11").Value = "54.14"
$ws.Range("E11").Value = "  +1.15%  "
$ws.Range("D12").Value = "0.0000306"
$ws.Range("E12").Value = "  +1.79%  "
$ws.Range("E13").Value = "  +0.80%  "
$ws.Range("D14").Value = "4.216.00"
$ws.Range("D15").Value = "678.63"
$ws.Range("E15").Value = "  +14.01%  "
$ws.Range("D16").Value = "13.04"
$ws.Range("E16").Value = "  +1.81%  "
$ws.Range("D17").Value = "71.065.40"
$ws.Range("E17").Value = "  +1.76%  "
$ws.Range("D18").Value = "3.638.76"
$ws.Range("E18").Value = "  +3.86%  "
$ws.Range("D19").Value = "19.11"
$ws.Range("E19").Value = "  -0.57%  "
$ws.Range("E20").Value = "  +0.58%  "
$ws.Range("E21").Value = "  +1.24%  "
$ws.Range("D22").Value = "18.47"
$ws.Range("E22").Value = "  +0.77%  "
$ws.Range("D23").Value = "5.37"
$ws.Range("E23").Value = "  +1.51%  "
$ws.Range("D24").Value = "106.34"
$ws.Range("E24").Value = "  +4.29%  "
$ws.Range("D25").Value = "4.63"
$ws.Range("E25").Value = "  -0.52%  "
$ws.Range("D26").Value = "3.00"
$ws.Range("E26").Value = "  -5.09%  "
$ws.Range("D27").Value = "10.46"
$ws.Range("E27").Value = "  -3.47%  "
$ws.Range("D28").Value = "9.81"
$ws.Range("E28").Value = "  +2.83%  "
$ws.Range("E29").Value = "  +2.71%  "
$ws.Range("D30").Value = "4.71"
$ws.Range("E30").Value = "  +9.81%  "
$ws.Range("D31").Value = "7.24"
$ws.Range("E31").Value = "  +2.46%  "
$ws.Range("D32").Value = "12.26"
$ws.Range("E32").Value = "  -1.33%  "
$ws.Range("E33").Value = "  +0.67%  "
$ws.Range("D34").Value = "63.54"
$ws.Range("E34").Value = "  +0.68%  "
$ws.Range("D35").Value = "3.982.55"
$ws.Range("E35").Value = "  +6.84%  "
$ws.Range("D36").Value = "0.0₃0874"
$ws.Range("E36").Value = "  +5.78%  "
$ws.Range("E37").Value = "  -0.16%  "
$ws.Range("D38").Value = "3.03"
$ws.Range("E38").Value = "  -1.93%  "
$ws.Range("D39").Value = "37.09"
$ws.Range("E39").Value = "  +1.87%  "
$ws.Range("D40").Value = "508.87"
$ws.Range("E40").Value = "  +5.37%  "
$ws.Range("E41").Value = "  -0.82%  "
$ws.Range("D42").Value = "3.55"
$ws.Range("E42").Value = "  -2.95%  "
$ws.Range("E43").Value = "  +2.31%  "
$ws.Range("D44").Value = "3.12"
$ws.Range("E44").Value = "  +10.52%  "
$ws.Range("D45").Value = "0.0462"
$ws.Range("D46").Value = "3.51"
$ws.Range("E46").Value = "  +6.97%  "
$ws.Range("E47").Value = "  +0.82%  "
$ws.Range("D48").Value = "8.71"
$ws.Range("E48").Value = "  +3.43%  "
$ws.Range("E49").Value = "  -0.23%  "
$ws.Range("E50").Value = "  +1.32%  "
$ws.Range("E51").Value = "  +4.82%  "

# Restore default style on column D so no extra number-format style lingers
# on cells (matches the original workbook, which has no explicit style there).
$ws.Range("D2:D51").Style = "Normal"

